$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B width (narrow "Step" column -> wider "Datatype" column) ---
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# --- Update existing formula text cells (C10:C14) to use new X/Y/name based sample ---
# Leading "'" forces these to be stored as literal text (quotePrefix) instead of formulas,
# matching the original workbook's convention for these cells.
$ws.Range("C10").Value = "'= new Y[] { new Y(`"a1`") } "
$ws.Range("C11").Value = "'= `$S1[(X x) transform to x.name]"
$ws.Range("C12").Value = "'= `$S1[(X x) transform unique to x.name]"
$ws.Range("C13").Value = "'= `$S1[(X x) select all having x.name.length > 0]"
$ws.Range("C14").Value = "'= `$S1[(X x) split by x.name.length > 0]"

# --- Add new "Datatype X" / "Datatype Y" definition blocks ---
# Detail rows (22, 26) first, so the centered/wrapped style is created once and reused.
$ws.Range("B22").Value = "String"
$ws.Range("C22").Value = "name"
$r22 = $ws.Range("B22:C22")
$r22.HorizontalAlignment = -4108
$r22.VerticalAlignment = -4108
$r22.WrapText = $true

$ws.Range("B26").Value = "String"
$ws.Range("C26").Value = "name"
$r26 = $ws.Range("B26:C26")
$r26.HorizontalAlignment = -4108
$r26.VerticalAlignment = -4108
$r26.WrapText = $true

# Header rows (21, 25), merged, same centered/wrapped style.
$ws.Range("B21").Value = "Datatype X"
$r21 = $ws.Range("B21:C21")
$r21.HorizontalAlignment = -4108
$r21.VerticalAlignment = -4108
$r21.WrapText = $true
$r21.MergeCells = $true

$ws.Range("B25").Value = "Datatype Y"
$r25 = $ws.Range("B25:C25")
$r25.HorizontalAlignment = -4108
$r25.VerticalAlignment = -4108
$r25.WrapText = $true
$r25.MergeCells = $true

# --- Update selection / scroll position ---
$ws.Range("F7").Select()
